$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) "Version" sheet: bump version number + last-updated date
# -----------------------------------------------------------------
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("B3").Value = 4.03
$wsVersion.Range("B4").Value = 45358

# -----------------------------------------------------------------
# 2) "Introduction" sheet: move the (stale) selection
# -----------------------------------------------------------------
$wsIntro = $wb.Worksheets.Item("Introduction")
$wsIntro.Activate()
$wsIntro.Range("B54").Select()

# -----------------------------------------------------------------
# 3) "Skills" sheet: hide the helper columns D:H and move selection
# -----------------------------------------------------------------
$wsSkills = $wb.Worksheets.Item("Skills")
$wsSkills.Activate()
$wsSkills.Range("D1:H1048576").Select()
$wsSkills.Columns("D:G").Hidden = $true
$wsSkills.Columns("H").ColumnWidth = 0
$wsSkills.Columns("H").Hidden = $true

# -----------------------------------------------------------------
# 4) "Goal Outline" sheet: add a third year block (2026) mirroring
#    the existing 2024/2025 layout, shifted 12 columns to the right
#    (N:Y -> Z:AK)
# -----------------------------------------------------------------
$wsGoal = $wb.Worksheets.Item("Goal Outline")
$wsGoal.Activate()

# Header band (year banner row) + spacer row - copy as full blocks
# since every cell in N5:Y6 already exists (no gaps).
$wsGoal.Range("N5:Y5").Copy($wsGoal.Range("Z5"))
$wsGoal.Range("N6:Y6").Copy($wsGoal.Range("Z6"))

# Body rows 7-20 only have a formatted cell in column N (mirrored to
# Z) - copy one cell at a time so we don't materialise the empty
# gap cells (O..Y) that a whole-range copy would stamp out.
for ($r = 7; $r -le 20; $r++) {
    $wsGoal.Range("N$r").Copy($wsGoal.Range("Z$r"))
}

# Month header row - full block again, no gaps.
$wsGoal.Range("N21:Y21").Copy($wsGoal.Range("Z21"))

# Fix up the year values: shift existing years forward one column
# group and set the new block to 2026.
$wsGoal.Range("B5:M5").Value = 2024
$wsGoal.Range("N5:Y5").Value = 2025
$wsGoal.Range("Z5").Value = 2026

# Recreate the merged header cell for the new block.
$wsGoal.Range("Z5:AK5").Merge()

$wsGoal.Range("Z6").Select()

# -----------------------------------------------------------------
# Restore the originally active sheet/tab.
# -----------------------------------------------------------------
$wsVersion.Activate()
